# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values on Sheet1 for rows 2-24 (excluding row 20, unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 2
    4  = 2
    5  = 3
    6  = 0
    7  = 1
    8  = 2
    9  = 0
    10 = 2
    11 = 2
    12 = 1
    13 = 3
    14 = 3
    15 = 2
    16 = 1
    17 = 0
    18 = 1
    19 = 1
    21 = 1
    22 = 1
    23 = 2
    24 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
